$wb = $excel.ActiveWorkbook

# New week header label to add in column AN (follows the existing weekly date-range
# headers in row 2, one week after "16.09-22.09.19" which lives in column AM).
# NOTE: we read/write through .Value2 rather than .Value - in this host the .Value
# getter resolves to a member-descriptor instead of invoking it, so .Value2 (which
# behaves correctly both ways) is used everywhere for reliability.
$newWeekHeader = "23.09-29.09.19"

$ws = $wb.Worksheets.Item("Kundenjobs")
$ws.Range("AN2").Value2 = $newWeekHeader
# Column AN was left blank for these rows even though AM (the prior week) has a
# value; carry that same figure forward into AN so the new week column is populated.
$ws.Range("AN8").Value2 = $ws.Range("AM8").Value2
$ws.Range("AN9").Value2 = $ws.Range("AM9").Value2
$ws.Range("AN12").Value2 = $ws.Range("AM12").Value2
$ws.Range("AN13").Value2 = $ws.Range("AM13").Value2
$ws.Range("AN14").Value2 = $ws.Range("AM14").Value2
$ws.Range("AN15").Value2 = $ws.Range("AM15").Value2
$ws.Range("AN16").Value2 = $ws.Range("AM16").Value2
$ws.Range("AN17").Value2 = $ws.Range("AM17").Value2
$ws.Range("AN19").Value2 = $ws.Range("AM19").Value2
$ws.Range("AN20").Value2 = $ws.Range("AM20").Value2
$ws.Range("AN21").Value2 = $ws.Range("AM21").Value2
$ws.Range("AN22").Value2 = $ws.Range("AM22").Value2
$ws.Range("AN25").Value2 = $ws.Range("AM25").Value2
$ws.Range("AN26").Value2 = $ws.Range("AM26").Value2
$ws.Range("AN27").Value2 = $ws.Range("AM27").Value2
$ws.Range("AN28").Value2 = $ws.Range("AM28").Value2
$ws.Range("AN29").Value2 = $ws.Range("AM29").Value2
$ws.Range("AN30").Value2 = $ws.Range("AM30").Value2
$ws.Range("AN32").Value2 = $ws.Range("AM32").Value2
$ws.Range("AN33").Value2 = $ws.Range("AM33").Value2
$ws.Range("AN34").Value2 = $ws.Range("AM34").Value2
$ws.Range("AN35").Value2 = $ws.Range("AM35").Value2
$ws.Range("AN41").Value2 = $ws.Range("AM41").Value2
$ws.Range("AN42").Value2 = $ws.Range("AM42").Value2
$ws.Range("AN44").Value2 = $ws.Range("AM44").Value2
$ws.Range("AN46").Value2 = $ws.Range("AM46").Value2
$ws.Range("AN51").Value2 = $ws.Range("AM51").Value2
$ws.Range("AN52").Value2 = $ws.Range("AM52").Value2
$ws.Range("AN56").Value2 = $ws.Range("AM56").Value2
$ws.Range("AN60").Value2 = $ws.Range("AM60").Value2
$ws.Range("AN61").Value2 = $ws.Range("AM61").Value2
$ws.Range("AN62").Value2 = $ws.Range("AM62").Value2
$ws.Range("AN64").Value2 = $ws.Range("AM64").Value2
$ws.Range("AN65").Value2 = $ws.Range("AM65").Value2

$ws = $wb.Worksheets.Item("Pitch_Neugeschäft")
$ws.Range("AN2").Value2 = $newWeekHeader

$ws = $wb.Worksheets.Item("Keine Arbeit")
$ws.Range("AN2").Value2 = $newWeekHeader
# Column AN was left blank for these rows even though AM (the prior week) has a
# value; carry that same figure forward into AN so the new week column is populated.
$ws.Range("AN14").Value2 = $ws.Range("AM14").Value2
$ws.Range("AN15").Value2 = $ws.Range("AM15").Value2
$ws.Range("AN16").Value2 = $ws.Range("AM16").Value2
$ws.Range("AN19").Value2 = $ws.Range("AM19").Value2
$ws.Range("AN20").Value2 = $ws.Range("AM20").Value2
$ws.Range("AN22").Value2 = $ws.Range("AM22").Value2
$ws.Range("AN27").Value2 = $ws.Range("AM27").Value2
$ws.Range("AN30").Value2 = $ws.Range("AM30").Value2
$ws.Range("AN32").Value2 = $ws.Range("AM32").Value2
$ws.Range("AN33").Value2 = $ws.Range("AM33").Value2
$ws.Range("AN42").Value2 = $ws.Range("AM42").Value2
$ws.Range("AN44").Value2 = $ws.Range("AM44").Value2
$ws.Range("AN52").Value2 = $ws.Range("AM52").Value2
$ws.Range("AN62").Value2 = $ws.Range("AM62").Value2
$ws.Range("AN64").Value2 = $ws.Range("AM64").Value2
$ws.Range("AN65").Value2 = $ws.Range("AM65").Value2

$ws = $wb.Worksheets.Item("Interne Jobs")
$ws.Range("AN2").Value2 = $newWeekHeader
# Column AN was left blank for these rows even though AM (the prior week) has a
# value; carry that same figure forward into AN so the new week column is populated.
$ws.Range("AN8").Value2 = $ws.Range("AM8").Value2
$ws.Range("AN9").Value2 = $ws.Range("AM9").Value2
$ws.Range("AN12").Value2 = $ws.Range("AM12").Value2
$ws.Range("AN13").Value2 = $ws.Range("AM13").Value2
$ws.Range("AN14").Value2 = $ws.Range("AM14").Value2
$ws.Range("AN15").Value2 = $ws.Range("AM15").Value2
$ws.Range("AN16").Value2 = $ws.Range("AM16").Value2
$ws.Range("AN17").Value2 = $ws.Range("AM17").Value2
$ws.Range("AN19").Value2 = $ws.Range("AM19").Value2
$ws.Range("AN20").Value2 = $ws.Range("AM20").Value2
$ws.Range("AN21").Value2 = $ws.Range("AM21").Value2
$ws.Range("AN22").Value2 = $ws.Range("AM22").Value2
$ws.Range("AN25").Value2 = $ws.Range("AM25").Value2
$ws.Range("AN26").Value2 = $ws.Range("AM26").Value2
$ws.Range("AN27").Value2 = $ws.Range("AM27").Value2
$ws.Range("AN28").Value2 = $ws.Range("AM28").Value2
$ws.Range("AN29").Value2 = $ws.Range("AM29").Value2
$ws.Range("AN30").Value2 = $ws.Range("AM30").Value2
$ws.Range("AN32").Value2 = $ws.Range("AM32").Value2
$ws.Range("AN33").Value2 = $ws.Range("AM33").Value2
$ws.Range("AN34").Value2 = $ws.Range("AM34").Value2
$ws.Range("AN35").Value2 = $ws.Range("AM35").Value2
$ws.Range("AN41").Value2 = $ws.Range("AM41").Value2
$ws.Range("AN42").Value2 = $ws.Range("AM42").Value2
$ws.Range("AN44").Value2 = $ws.Range("AM44").Value2
$ws.Range("AN46").Value2 = $ws.Range("AM46").Value2
$ws.Range("AN51").Value2 = $ws.Range("AM51").Value2
$ws.Range("AN52").Value2 = $ws.Range("AM52").Value2
$ws.Range("AN53").Value2 = $ws.Range("AM53").Value2
$ws.Range("AN56").Value2 = $ws.Range("AM56").Value2
$ws.Range("AN60").Value2 = $ws.Range("AM60").Value2
$ws.Range("AN61").Value2 = $ws.Range("AM61").Value2
$ws.Range("AN62").Value2 = $ws.Range("AM62").Value2
$ws.Range("AN64").Value2 = $ws.Range("AM64").Value2
$ws.Range("AN65").Value2 = $ws.Range("AM65").Value2

$ws = $wb.Worksheets.Item("Urlaub")
$ws.Range("AN2").Value2 = $newWeekHeader
# Column AN was left blank for these rows even though AM (the prior week) has a
# value; carry that same figure forward into AN so the new week column is populated.
$ws.Range("AN13").Value2 = $ws.Range("AM13").Value2
$ws.Range("AN19").Value2 = $ws.Range("AM19").Value2
$ws.Range("AN20").Value2 = $ws.Range("AM20").Value2
$ws.Range("AN21").Value2 = $ws.Range("AM21").Value2
$ws.Range("AN29").Value2 = $ws.Range("AM29").Value2
$ws.Range("AN31").Value2 = $ws.Range("AM31").Value2
$ws.Range("AN33").Value2 = $ws.Range("AM33").Value2
$ws.Range("AN61").Value2 = $ws.Range("AM61").Value2
$ws.Range("AN62").Value2 = $ws.Range("AM62").Value2

$ws = $wb.Worksheets.Item("Krankheit")
$ws.Range("AN2").Value2 = $newWeekHeader
# Column AN was left blank for these rows even though AM (the prior week) has a
# value; carry that same figure forward into AN so the new week column is populated.
$ws.Range("AN22").Value2 = $ws.Range("AM22").Value2
$ws.Range("AN26").Value2 = $ws.Range("AM26").Value2
$ws.Range("AN65").Value2 = $ws.Range("AM65").Value2

$ws = $wb.Worksheets.Item("Feiertage")
$ws.Range("AN2").Value2 = $newWeekHeader

$ws = $wb.Worksheets.Item("Überstundenabbau")
$ws.Range("AN2").Value2 = $newWeekHeader
